$d = $word.ActiveDocument
$d.Content.Find.Execute("This is a Microsoft word document.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "This is a Microsoft word document. (Changed main)", 2)
